$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D held boolean formulas (=TRUE()/=FALSE()) that evaluated to 1/0.
# The fix stores plain text "TRUE"/"FALSE" instead. Assigning the bare word
# "TRUE"/"FALSE" to a cell's .Value auto-converts it to a real boolean, so
# the literal text is produced on scratch cells via a formula and brought
# into the target cells with a values-only paste (keeps the text type and
# the column's existing "@" text style instead of a formula or a number).
$scratchTrue = $ws.Range("Z1")
$scratchFalse = $ws.Range("Z2")
$scratchTrue.Formula = '="TRUE"'
$scratchFalse.Formula = '="FALSE"'

$scratchTrue.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").PasteSpecial(-4163)

$scratchFalse.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("D7").PasteSpecial(-4163)

$scratchTrue.ClearContents()
$scratchFalse.ClearContents()

# Update the stored selection to match the authored workbook state.
$ws.Range("E7").Select()
